{"js": "// The document repeatedly referred to \"linear regression\" when the labs\n// are actually about the \"logistic regression\" model used for the\n// diabetes / bike-rental classification experiments. Replace every\n// whole-word, case-sensitive occurrence of \"linear\" -> \"logistic\" and\n// \"Linear\" -> \"Logistic\" throughout the document body.\n\nconst body = context.document.body;\n\nconst lower = body.search(\"linear\", { matchCase: true, matchWholeWord: true });\nlower.load(\"items\");\nawait context.sync();\n\nfor (const range of lower.items) {\n  range.insertText(\"logistic\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\nconst upper = body.search(\"Linear\", { matchCase: true, matchWholeWord: true });\nupper.load(\"items\");\nawait context.sync();\n\nfor (const range of upper.items) {\n  range.insertText(\"Logistic\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The document used \"linear regression\" terminology where it should have\n# read \"logistic regression\" (the labs are about classification models).\n# Replace every whole-word occurrence of \"linear\"/\"Linear\" with\n# \"logistic\"/\"Logistic\" throughout the document body, preserving case.\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"linear\", $true, $true, $false, $false, $false, $true, 1, $false, \"logistic\", 2) | Out-Null\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\"Linear\", $true, $true, $false, $false, $false, $true, 1, $false, \"Logistic\", 2) | Out-Null\n"}
